# This edit removes the "Elab", "X", "Q2" and "target" columns from the
# data table (originally columns A, C, D and F), shifting the remaining
# columns left. The stray "N/A" value that lands in column B after the
# shift is corrected to "SU3" (its real value), and the selection /
# active cell is moved to E7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the unwanted columns, working from right to left so earlier
# deletions don't shift the column letters of the ones still to be
# removed.
$ws.Range("F:F").EntireColumn.Delete()   # was "target"
$ws.Range("D:D").EntireColumn.Delete()   # was "Q2"
$ws.Range("C:C").EntireColumn.Delete()   # was "X"
$ws.Range("A:A").EntireColumn.Delete()   # was "Elab"

# After the shift, B2 (formerly E2, "N/A") needs to read "SU3" to match
# what was originally column F's value.
$ws.Range("B2").Value = "SU3"

# Move the selection to match the saved view state.
$ws.Range("E7").Select()
